$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1316.7059
$ws.Range("I53").Value = 162.6
$ws.Range("J53").Value = 1797.5834
$ws.Range("K53").Value = 162.6
$ws.Range("L53").Value = 1797.5834
$ws.Range("M53").Value = 474.4
$ws.Range("N53").Value = -3071.5834

$ws.Range("H64").Value = 6206.5
$ws.Range("J64").Value = 6997
$ws.Range("L64").Value = 6997
$ws.Range("N64").Value = -7493

$ws.Range("H67").Value = 6206.5
$ws.Range("J67").Value = 6997
$ws.Range("L67").Value = 6997
$ws.Range("N67").Value = -8713

$ws.Range("H80").Value = 3068.5
$ws.Range("J80").Value = 6098.6
$ws.Range("L80").Value = 18295.8
$ws.Range("N80").Value = -20291.8

$ws.Range("H83").Value = 3068.5
$ws.Range("J83").Value = 6098.6
$ws.Range("L83").Value = 54887.4
$ws.Range("N83").Value = -64871.4

$ws.Range("H88").Value = 13914.5
$ws.Range("I88").Value = 3996
$ws.Range("J88").Value = 17220.666
$ws.Range("K88").Value = 3996
$ws.Range("L88").Value = 17220.666
$ws.Range("M88").Value = -3590
$ws.Range("N88").Value = -18032.666

$ws.Range("H91").Value = 13914.5
$ws.Range("I91").Value = 3996
$ws.Range("J91").Value = 17220.666
$ws.Range("K91").Value = 3996
$ws.Range("L91").Value = 17220.666
$ws.Range("M91").Value = -2592
$ws.Range("N91").Value = -20028.666

$ws.Range("H135").Value = 1252.8823
$ws.Range("I135").Value = 1019.86664
$ws.Range("K135").Value = 9178.79976
$ws.Range("M135").Value = -6643.79976

$ws.Range("H138").Value = 3536.02
$ws.Range("I138").Value = 2852.25
$ws.Range("J138").Value = 3857.7942
$ws.Range("K138").Value = 8556.75
$ws.Range("L138").Value = 11573.3826
$ws.Range("M138").Value = -3416.75
$ws.Range("N138").Value = -21853.3826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 5000
$ws.Range("K31").Value = 5000
$ws.Range("M31").Value = -4706

$ws.Range("H32").Value = 1646.1177
$ws.Range("I32").Value = 1250.3438
$ws.Range("J32").Value = 7978.5
$ws.Range("K32").Value = 1250.3438
$ws.Range("L32").Value = 7978.5
$ws.Range("M32").Value = -963.3438000000001
$ws.Range("N32").Value = -8552.5

$ws.Range("H74").Value = 15152778
$ws.Range("I74").Value = 19608890
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 19608890
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -19608016
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 15152778
$ws.Range("I77").Value = 19608890
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 98044450
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -98040082
$ws.Range("N77").Value = -18736

$ws.Range("H92").Value = 26664
$ws.Range("J92").Value = 26664
$ws.Range("L92").Value = 26664
$ws.Range("N92").Value = -31656

$ws.Range("H122").Value = 2246.077
$ws.Range("J122").Value = 2678.2222
$ws.Range("L122").Value = 8034.6666
$ws.Range("N122").Value = -12934.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 3272.9
$ws.Range("I37").Value = 3466.6667
$ws.Range("J37").Value = 2982.25
$ws.Range("K37").Value = 3466.6667
$ws.Range("L37").Value = 2982.25
$ws.Range("M37").Value = -3329.6667
$ws.Range("N37").Value = -3256.25

$ws.Range("H86").Value = 3989
$ws.Range("I86").Value = 1848.75
$ws.Range("J86").Value = 5701.2
$ws.Range("K86").Value = 1848.75
$ws.Range("L86").Value = 5701.2
$ws.Range("M86").Value = -725.75
$ws.Range("N86").Value = -7947.2

$ws.Range("H89").Value = 3989
$ws.Range("I89").Value = 1848.75
$ws.Range("J89").Value = 5701.2
$ws.Range("K89").Value = 9243.75
$ws.Range("L89").Value = 28506
$ws.Range("M89").Value = -3627.75
$ws.Range("N89").Value = -39738

$ws.Range("H134").Value = 1391.0526
$ws.Range("I134").Value = 1381.9412
$ws.Range("K134").Value = 4145.8236
$ws.Range("M134").Value = -1610.8236

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 39995
$ws.Range("I82").Value = 39995
$ws.Range("K82").Value = 39995
$ws.Range("M82").Value = -39634

$ws.Range("H85").Value = 39995
$ws.Range("I85").Value = 39995
$ws.Range("K85").Value = 39995
$ws.Range("M85").Value = -38747

$ws.Range("H107").Value = 803.86957
$ws.Range("I107").Value = 605.86664
$ws.Range("J107").Value = 1175.125
$ws.Range("K107").Value = 605.86664
$ws.Range("L107").Value = 1175.125
$ws.Range("M107").Value = 1314.13336
$ws.Range("N107").Value = -5015.125

$ws.Range("H122").Value = 5655.077
$ws.Range("I122").Value = 1865.7142
$ws.Range("J122").Value = 10076
$ws.Range("K122").Value = 5597.142599999999
$ws.Range("L122").Value = 30228
$ws.Range("M122").Value = -3147.142599999999
$ws.Range("N122").Value = -35128

$ws.Range("H132").Value = 16503.25
$ws.Range("I132").Value = 15666.333
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 46998.999
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -44468.999
$ws.Range("N132").Value = -62102

$ws.Range("H134").Value = 7484.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 7484.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 22453.5
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -27523.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4565.143
$ws.Range("I132").Value = 5073.75
$ws.Range("J132").Value = 4361.7
$ws.Range("K132").Value = 45663.75
$ws.Range("L132").Value = 39255.3
$ws.Range("M132").Value = -43133.75
$ws.Range("N132").Value = -44315.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 31000
$ws.Range("J117").Value = 31000
$ws.Range("L117").Value = 31000
$ws.Range("N117").Value = -37884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6356.857
$ws.Range("I82").Value = 6712.625
$ws.Range("J82").Value = 5882.5
$ws.Range("K82").Value = 6712.625
$ws.Range("L82").Value = 5882.5
$ws.Range("M82").Value = -6351.625
$ws.Range("N82").Value = -6604.5

$ws.Range("H85").Value = 6356.857
$ws.Range("I85").Value = 6712.625
$ws.Range("J85").Value = 5882.5
$ws.Range("K85").Value = 6712.625
$ws.Range("L85").Value = 5882.5
$ws.Range("M85").Value = -5464.625
$ws.Range("N85").Value = -8378.5

$ws.Range("H132").Value = 4658.3667
$ws.Range("I132").Value = 2589.1738
$ws.Range("J132").Value = 11457.143
$ws.Range("K132").Value = 7767.5214
$ws.Range("L132").Value = 34371.429
$ws.Range("M132").Value = -5237.5214
$ws.Range("N132").Value = -39431.429

$ws.Range("H140").Value = 69895
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 69895
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 69895
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -80255

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1380.4
$ws.Range("I100").Value = 671.1429000000001
$ws.Range("J100").Value = 3035.3333
$ws.Range("K100").Value = 1342.2858
$ws.Range("L100").Value = 6070.6666
$ws.Range("M100").Value = -801.2858000000001
$ws.Range("N100").Value = -7152.6666

$ws.Range("H126").Value = 2572.2666
$ws.Range("I126").Value = 1330.3636
$ws.Range("J126").Value = 5987.5
$ws.Range("K126").Value = 3991.0908
$ws.Range("L126").Value = 17962.5
$ws.Range("M126").Value = -1521.0908
$ws.Range("N126").Value = -22902.5

$ws.Range("H132").Value = 7414.875
$ws.Range("I132").Value = 6642.2
$ws.Range("K132").Value = 19926.6
$ws.Range("M132").Value = -17396.6
